$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a brand-new entry ("Adenomyoma") as the new row 3. This pushes every
# following data row down by one (their values/styles move with them).
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = "Gallbladder and biliary tract"
$ws.Cells.Item(3, 2).Value = "Adenomyoma "
$ws.Cells.Item(3, 3).Value = "Clip 1 B-mode + color + microV"
$ws.Cells.Item(3, 4).Value = "https://youtu.be/ZXwd0gwHEkQ "

# ---------------------------------------------------------------------------
# The row insert does not carry the worksheet's hyperlinks along with it, so
# rebuild the whole hyperlink collection against the (now shifted) D column.
# ---------------------------------------------------------------------------
$ws.Range("D3").Hyperlinks.Delete()

$links = @{
    "D4"  = "https://youtu.be/zxTC0YBY2RY"
    "D5"  = "https://youtu.be/K2Wbg7BgXy4"
    "D7"  = "https://youtu.be/2kRZcpi70Aw"
    "D10" = "https://youtu.be/91M82AIMyu0"
    "D13" = "https://youtu.be/15o_Km86IzM"
    "D15" = "https://youtu.be/RhSUFLTmTl4"
    "D19" = "https://youtu.be/DjI1kEnzfSQ"
    "D20" = "https://youtu.be/U3ydTsRwxok"
    "D25" = "https://youtu.be/xBfd04F4Ni8"
    "D26" = "https://youtu.be/JvwODCASLYQ"
    "D27" = "https://youtu.be/pc-vbxSRTbs"
    "D28" = "https://youtu.be/Axbee4vjNtU"
    "D29" = "https://youtu.be/qushjTAy6XQ"
    "D30" = "https://youtu.be/_FckFwJwynI"
    "D31" = "https://youtu.be/z_oaRVxRz5s"
    "D3"  = "https://youtu.be/ZXwd0gwHEkQ"
}

foreach ($addr in $links.Keys) {
    $ws.Hyperlinks.Add($ws.Range($addr), $links[$addr])
    $ws.Range($addr).Style = "Collegamento ipertestuale"
}

# ---------------------------------------------------------------------------
# Refresh the saved sort-state range so it also covers the newly added row.
# ---------------------------------------------------------------------------
$sortObj = $ws.Sort
$sortObj.SetRange($ws.Range("A2:C29"))
$sortObj.Header = 0
$sortObj.Apply()

$ws.Range("D6").Select()
